# Update database: drop the oldest reporting period (column D, "6 ماهه منتهی
# به 1399/06") and append the newest one (new last column, "12 ماهه منتهی به
# 1401/12") together with its data, per the updated read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest period column (D). This shifts every later column one
#    slot to the left (E->D, F->E, ... M->L) and cleans up the now-unused
#    shared string for the removed header/label.
$ws.Columns("D").Delete()

# 2) Clone the formatting of the (new) last data column L into the new last
#    column M so the freshly appended column inherits the same styles
#    (fills/borders/alignment) as the rest of the table.
$ws.Range("L8:L27").Copy()
$ws.Range("M8:M27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new column M with the newest reporting period.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-23"

$ws.Range("M11").Value = 106811
$ws.Range("M12").Value = -74576
$ws.Range("M13").Value = 32234
$ws.Range("M14").Value = -9452
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = -120
$ws.Range("M17").Value = 22663
$ws.Range("M18").Value = -7061
$ws.Range("M19").Value = 3266
$ws.Range("M20").Value = 18868
$ws.Range("M21").Value = -147
$ws.Range("M22").Value = 18721
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 18721
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 11335
$ws.Range("M27").Value = 0

# 4) The 9-month-1401 period's publish-date annotation is amended to
#    reflect the later re-publish (now the "(8)"-th revision instead of the
#    "(7)"-th, dated 1402-02-23).
$ws.Range("I9").Value = "1402-02-23 (8)"

# 5) Column M needs the same "wide" column width as the other
#    rightmost/odd columns (28 vs 29 alternating pattern).
$ws.Columns("M").ColumnWidth = 29
